$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the now-unused last row (row 10) ---------------------------
# The table shrinks from 10 data rows to 9; row 10 is dropped entirely.
$ws.Rows("10").Delete()

# --- Add the two new columns: Village (L) and Agronomist (M) -----------

# Header row (row 1)
$ws.Range("L1").Value = "Village"
$ws.Range("M1").Value = "Agronomist"
$ws.Range("K1").Copy()
$ws.Range("L1:M1").PasteSpecial(-4122)   # xlPasteFormats - match header styling

# Data row (row 2)
$ws.Range("L2").Value = "Ukpo"
$ws.Range("M2").Value = "Paul walker"
$ws.Range("K2").Copy()
$ws.Range("L2:M2").PasteSpecial(-4122)   # xlPasteFormats - match data row styling

# Blank filler rows (rows 3-9) - just extend the same look as column K
$ws.Range("K3").Copy()
$ws.Range("L3:M9").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

Write-Host "Activities template updated"
